# Apply the "Add data for 2021-09-28" update to the carjacking-by-neighborhood-by-month workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet to reflect the new "through" date.
$ws.Name = "Through 2021-09-20"

# Update the September 2021 column header text.
$ws.Range("B1").Value = "September 2021 (through September 20)"

# Update existing cell values (incremented counts from the newly added incident data).
$ws.Range("B3").Value = 7
$ws.Range("T3").Value = 5
$ws.Range("K4").Value = 6
$ws.Range("K5").Value = 9
$ws.Range("K6").Value = 3
$ws.Range("B7").Value = 7
$ws.Range("K7").Value = 2
$ws.Range("B11").Value = 3
$ws.Range("B17").Value = 3
$ws.Range("K17").Value = 4
$ws.Range("B23").Value = 3
$ws.Range("B28").Value = 3
$ws.Range("B55").Value = 4
$ws.Range("K55").Value = 7
$ws.Range("BD55").Value = 2
$ws.Range("B56").Value = 3
$ws.Range("AL64").Value = 2
$ws.Range("B74").Value = 2
$ws.Range("AC99").Value = 2

# New cells that did not previously have a value.
$ws.Range("B8").Value = 1
$ws.Range("AC42").Value = 1
$ws.Range("B65").Value = 1
